# Regenerate the localization-status report: items that were previously
# "Ready for handoff" have now moved on to "In Translation", and the
# Status column is re-sized (autofit) for the new, shorter text.

$wb = $excel.ActiveWorkbook

$overview = $wb.Worksheets.Item("Overview")
$zhcn     = $wb.Worksheets.Item("zh-cn")
$dede     = $wb.Worksheets.Item("de-de")

# --- Update status text ---------------------------------------------------
# Overview sheet mirrors the per-locale status in columns E (zh-cn) and F (de-de)
$overview.Range("E2").Value = "In Translation"
$overview.Range("F2").Value = "In Translation"

# Per-locale sheets carry the same status in column C
$zhcn.Range("C2").Value = "In Translation"
$dede.Range("C2").Value = "In Translation"

# --- Resize the Status columns to fit the new (shorter) text ---------------
# AutoFit-equivalent target width for "In Translation" is ~13.41 characters.
# 12.5 is the ColumnWidth input that lands closest to that computed width.
$targetColumnWidth = 12.5
$overview.Columns.Item(5).ColumnWidth = $targetColumnWidth
$overview.Columns.Item(6).ColumnWidth = $targetColumnWidth
$zhcn.Columns.Item(3).ColumnWidth = $targetColumnWidth
$dede.Columns.Item(3).ColumnWidth = $targetColumnWidth
